$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Market")

# Type the "market.other.menu" entry first (so it lands in the shared
# string table before the hardware entry below), right after the current
# "market.mod.menu" row (row 67), keeping column B sorted alphabetically.
$ws.Rows.Item(68).Insert()
$ws.Range("A68").Value = "cs"
$ws.Range("B68").Value = "market.other.menu"
$ws.Range("C68").Value = "Ostatní"

# Now add "market.hardware.menu" just above "market.home.menu" (row 56),
# keeping column B sorted alphabetically.
$ws.Rows.Item(56).Insert()
$ws.Range("A56").Value = "cs"
$ws.Range("B56").Value = "market.hardware.menu"
$ws.Range("C56").Value = "Hardware"

# Refresh the sheet's remembered sort range/state to cover the two new rows.
$sortRange = $ws.Range("A2:C89")
$keyRange = $ws.Range("B82:B89")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Restore the active cell/selection like the saved workbook had it.
$null = $ws.Range("B79").Select()
